# Merge NguyenTrungNghia with NguyenTrungNghia-updated
#
# The sheet holds two stacked "use case flow" tables:
#   rows 1-12  -> "View dock's information" flow
#   rows 14-24 -> "View detail information of bikes at the dock" flow
# The edit rewrites the text of the first table (dock info + available
# bikes) and the "Action" column of the second table (bike detail info),
# normalises every row height to 13.2, moves the active selection, and
# turns on an explicit (portrait) page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 1-12): "View dock's information" -----------------
$ws.Range("A1").Value = "View dock's information and avaiable bikes"
$ws.Range("C4").Value = "click on a dock to see dock's detail information and avaiable bikes"
$ws.Range("C5").Value = "update the number of avaiable bikes and e-bike's battery information"
$ws.Range("C6").Value = "display the view of dock's information and avaiable bikes"
$ws.Range("C10").Value = "in case of failed updating, notify to the user, end of use case"

# --- Table 2 (rows 14-24): "View detail information of bikes" -------
$ws.Range("C18").Value = "query for detail bike's information related to that bike"
$ws.Range("C19").Value = "display detail bike's information"
$ws.Range("C17").Value = "click on a bike in the list to see detail bike's information"

# --- Uniform row height across the whole used range ------------------
$ws.Range("A1:C24").RowHeight = 13.2

# --- Selection / scroll position --------------------------------------
$null = $ws.Range("C19").Select()

# --- Explicit page setup (portrait) -----------------------------------
$ws.PageSetup.Orientation = 1
